# Apply the workbook update:
#  1. Update the confidential-notice text in A18 (model-holdings-as-of date:
#     2021-04-09 -> 2021-04-21).
#  2. Update the Weight (D) and Percent Change (E) values for rows 2-15.
#
# The worksheet is protected, so it must be unprotected before writing and
# re-protected (same password) afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$password = "D382"
$ws.Unprotect($password)

# --- 1. Confidential notice text (row 18, column A) ---------------------
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-21 for illustrative purposes only and are subject to change."

# --- 2. Weight / Percent Change values (rows 2-15, columns D & E) -------
$ws.Range("D2").Value = 0.05687126241897092
$ws.Range("E2").Value = 0.01192982456140346

$ws.Range("D3").Value = 0.02345291042842685
$ws.Range("E3").Value = 0.009132420091324311

$ws.Range("D4").Value = 0.03129116722778537
$ws.Range("E4").Value = 0.007843887507174241

$ws.Range("D5").Value = 0.03034397801509982
$ws.Range("E5").Value = 0.01124531445231147

$ws.Range("D6").Value = 0.03574308955960025
$ws.Range("E6").Value = 0.01607860652076809

$ws.Range("D7").Value = 0.01865244375149453
$ws.Range("E7").Value = 0.01397903145282053

$ws.Range("D8").Value = 0.004594599358585281
$ws.Range("E8").Value = 0.01976112920738338

$ws.Range("D9").Value = 0.006833863743772014
$ws.Range("E9").Value = 0.0118746350009733

$ws.Range("D10").Value = 0.07052036539952394
$ws.Range("E10").Value = 0.009620826259196358

$ws.Range("D11").Value = 0.07060018471519969
$ws.Range("E11").Value = 0.009609949123798645

$ws.Range("D12").Value = 0.1484639271568925
$ws.Range("E12").Value = 0.002007168458781461

$ws.Range("D13").Value = 0.3884546681162741
$ws.Range("E13").Value = 0.0005239717055278703

$ws.Range("D14").Value = 0.1141775401083747
$ws.Range("E14").Value = 0.00477122583802303

$ws.Range("E15").Value = 0.004889930826585287

$ws.Protect($password)
